$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: mirror of row 2's data (per diff), appended below the existing data.
$ws.Range("A3").Value = "last"
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = "first"

# E3/F3 hold numeric-looking text ("3333"/"333"); force text so they are
# stored as strings rather than being coerced to numbers.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3333"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "333"

$ws.Range("G3").Value = "HCM"
$ws.Range("H3").Value = "address updated"

# Drop the text-format style picked up above so the new row matches row 2's
# formatting (no explicit style applied).
$ws.Range("E3:F3").Style = "Normal"
